$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: rename the "##" marker row labels ---
$ws.Range("A1").Value = "##var"
$ws.Range("A2").Value = "##type"
$ws.Range("A3").Value = "##"

# --- New rows 10-11 (SpawnConfig path entries re-using existing path names) ---
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = "路径1"
$ws.Range("D10").Value = "0;0;0"
$ws.Range("E10").Value = "1001;1002;1003;1004"

$ws.Range("B11").Value = 7
$ws.Range("C11").Value = "路径2"
$ws.Range("D11").Value = "0;0;10"
$ws.Range("E11").Value = 1006

# --- Sheet view: move the active selection cell ---
$ws.Range("J12").Select()
